$d = $word.ActiveDocument

# --- Step 1: fix "Baz chan" + bookmark + "ges" -> single clean run "Baz changes" ---
# Locate the paragraph that currently reads "Baz changes" (it is split across two
# runs with a _GoBack bookmark sandwiched between them in the source document).
$bazPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Baz changes" + [char]13) {
        $bazPara = $p
    }
}
$bazRange = $bazPara.Range
$bazRange.MoveEnd(1, -1) | Out-Null
$bazRange.Delete()
$bazRange.InsertAfter("Baz changes")

# --- Step 2: insert a brand-new paragraph right after it ---
$bazPara2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Baz changes" + [char]13) {
        $bazPara2 = $p
    }
}
$bazPara2.Range.InsertParagraphAfter()

# Find the (now empty) freshly inserted paragraph that follows it.
$newPara = $bazPara2.Next()
$newRange = $newPara.Range
$newRange.MoveEnd(1, -1) | Out-Null

# --- Step 3: populate the new paragraph with the formatted runs, proofErr
#             spell-check markers and the relocated _GoBack bookmark, using a
#             raw WordOpenXML fragment so run-splitting/formatting is exact ---
# Note: the final "<w:p/>" in the fragment merges into the (already existing)
# host paragraph that the collapsed $newRange sits inside, so to additionally
# produce one brand-new blank paragraph after our content paragraph we need
# *two* trailing "<w:p/>" elements here (first = new paragraph, second = merge
# target / the paragraph that was already there).
$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Wade </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Gava' + [char]0x2019 + 's</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> user ' + [char]0x2013 + ' sleepw4lker187) </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>changes:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Version control systems are software tools that help development teams manage changes to source code over time. Version control software keeps track of modifications to code so that if a mistake is made, developers can turn back the clock and compare earlier versions</w:t></w:r><w:r><w:t xml:space="preserve"> of the code. This then allows developers to fix the code whilst minimizing disruption to other team members. Hooray!</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($xmlFragment) | Out-Null
